$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("tpDictionary")
$ws2.Cells.Item(3,1).Value = "studyId"
$ws2.Cells.Item(4,1).Value = "studyArm"
$ws2.Cells.Item(5,1).Value = "subjectId"
$ws2.Cells.Item(8,1).Value = "outputPathId"
$ws2.Cells.Item(25,1).Value = "dose"
$ws2.Cells.Item(26,1).Value = "route"
$ws2.Cells.Item(6,1).Value = "individualId"
